$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mitochondrial")
$ws.Activate()
$ws.Rows.Item(206).Insert()
$ws.Cells.Item(206, 1).Value = "cytochrome c oxidase subunit I"
$ws.Cells.Item(206, 2).Value = "COI"
